# Update the two-digit multiplication practice sheet to the regenerated
# set of problems (output re-generated at commit 1c8df47).
#
# Most cells are simple 1:1 text substitutions. One row (row 10 of the
# table) additionally drops one cell ("74x51=") and gains a new one
# ("96x87=") at the end - since the cell COUNT in that row doesn't
# change (one removed, one added), we express the whole row as five
# positional cell-text assignments instead of trying to literally
# delete/insert a table cell.

$d = $word.ActiveDocument

function ReplaceProblem($old, $new) {
    $d.Content.Find.Execute($old, $true, $false, $false, $false, $false, `
                             $true, 1, $false, $new, 2) | Out-Null
}

# --- simple 1:1 replacements (rows 1, 5, 15, 20 and part of row 10) ---
ReplaceProblem "17×55=" "67×29="
ReplaceProblem "12×37=" "79×41="
ReplaceProblem "97×13=" "64×73="
ReplaceProblem "98×44=" "14×43="
ReplaceProblem "28×38=" "15×95="

ReplaceProblem "20×42=" "28×34="
ReplaceProblem "96×18=" "90×17="
ReplaceProblem "20×87=" "93×23="
ReplaceProblem "74×77=" "48×66="
ReplaceProblem "38×37=" "42×40="

ReplaceProblem "87×53=" "26×15="
ReplaceProblem "93×38=" "95×17="
ReplaceProblem "42×20=" "84×64="
ReplaceProblem "57×70=" "31×19="
ReplaceProblem "15×76=" "94×18="

ReplaceProblem "76×39=" "18×62="
ReplaceProblem "16×62=" "12×34="
ReplaceProblem "23×35=" "12×30="
ReplaceProblem "50×80=" "86×14="
ReplaceProblem "41×11=" "50×26="

# --- row 10: "74×51=" cell is deleted and a new "96×87=" cell is
# appended at the row's end; net cell count is unchanged, so rewrite
# the row's five cells positionally, left to right. ---
$row = $d.Tables(1).Rows(10)
$row.Cells(1).Range.Text = "78×23="
$row.Cells(2).Range.Text = "51×50="
$row.Cells(3).Range.Text = "46×13="
$row.Cells(4).Range.Text = "85×35="
$row.Cells(5).Range.Text = "96×87="
